$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188 (Excel COM semantics: existing row 188
# and everything below it shift down by one, so old rows 188-218 become
# rows 189-219).
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly price record.
$ws.Range("A188").Value = 7
$ws.Range("B188").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C188").Value = "Ñuble"
$ws.Range("D188").Value = 44504
$ws.Range("E188").Value = 16
$ws.Range("F188").Value = 100114001
$ws.Range("G188").Value = "Papa"
$ws.Range("H188").Value = "Patagonia"
$ws.Range("I188").Value = "1a (guarda)"
$ws.Range("J188").Value = 240
$ws.Range("K188").Value = 7000
$ws.Range("L188").Value = 8000
$ws.Range("M188").Value = 7500
$ws.Range("N188").Value = "$/saco 25 kilos"
$ws.Range("O188").Value = "Provincia de Diguillín"
$ws.Range("P188").Value = 300
$ws.Range("Q188").Value = 25
$ws.Range("R188").Value = "Hortaliza"

# Match the date-number format used by the rest of column D.
$ws.Range("D188").NumberFormat = $ws.Range("D189").NumberFormat
